# Update recovery-function parameter naming/values in comp_type_dmg_algo
# (sheet: "changed parameter names for recovery functions for general distribution")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comp_type_dmg_algo")

# Header rename: recovery_mean -> recovery_param1 ; recovery_std -> recovery_param2
$ws.Range("O1").Value = "recovery_param1"
$ws.Range("P1").Value = "recovery_param2"

# New recovery_function + recovery_param1 / recovery_param2 values per damage-state row.
# The distribution used switched from Normal (mean/std) to Rayleigh (param1/param2).
$rows = @(
    @{Row=2;  Param1=0; Param2=3},
    @{Row=3;  Param1=1; Param2=7},
    @{Row=4;  Param1=2; Param2=28},
    @{Row=5;  Param1=5; Param2=100},
    @{Row=6;  Param1=0; Param2=1},
    @{Row=7;  Param1=1; Param2=3},
    @{Row=8;  Param1=2; Param2=7},
    @{Row=9;  Param1=2; Param2=14},
    @{Row=10; Param1=0; Param2=1},
    @{Row=11; Param1=0; Param2=1},
    @{Row=12; Param1=0; Param2=1},
    @{Row=13; Param1=0; Param2=1},
    @{Row=14; Param1=0; Param2=1},
    @{Row=15; Param1=0; Param2=1},
    @{Row=16; Param1=0; Param2=1},
    @{Row=17; Param1=0; Param2=1},
    @{Row=18; Param1=0; Param2=1},
    @{Row=19; Param1=0; Param2=1},
    @{Row=20; Param1=0; Param2=1},
    @{Row=21; Param1=0; Param2=1}
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 14).Value = "Rayleigh"
    $ws.Cells.Item($rowNum, 15).Value = $r.Param1
    $ws.Cells.Item($rowNum, 16).Value = $r.Param2
}

# Reflect the view state captured in the saved file: frozen pane scrolled to K10,
# with O22 as the active selection on the bottom-right pane.
$ws.Activate()
$ws.Range("O22").Select()
